# ServicesAnalysis.xlsx — "Add files via upload"
#
# Re-applies the state the workbook was in when it was last saved:
#   1. Three rows in the "Sentinel Defined" column (H) flip from "No" to
#      "Yes" (Ciscoasa/row 5, Dionaea/row 9, ElasticPot/row 10).
#   2. The window had scrolled down (while the header row stayed frozen)
#      so row 6 is the first row under the frozen pane, with B12 as the
#      active/selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Flip "Sentinel Defined" to "Yes" for the three services ---
$ws.Range("H5").Value = "Yes"
$ws.Range("H9").Value = "Yes"
$ws.Range("H10").Value = "Yes"

# --- 2. Scroll the frozen view down to row 6 and select B12 ---
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B12").Select()
